# Corregimos el teléfono en la plantilla ret islr arsa.
#
# El marcador del teléfono de la compañía tenía espacios sobrantes dentro
# de las llaves: "{ companiaContabTelefono }". Lo dejamos sin espacios:
# "{companiaContabTelefono}". Ese marcador vive en la tabla "DATOS DEL
# AGENTE DE RETENCIÓN", fila "Teléfono:", segunda columna; localizamos la
# celda por contenido para no tocar el marcador parecido
# "{proveedorTelefono}" que aparece en otra tabla del documento.

$d = $word.ActiveDocument

$targetCell = $null
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
        for ($ci = 1; $ci -le $tbl.Columns.Count; $ci++) {
            try {
                $cell = $tbl.Cell($ri, $ci)
            } catch {
                continue
            }
            if ($cell.Range.Text -like "*companiaContabTelefono*") {
                $targetCell = $cell
            }
        }
    }
}

if ($targetCell -ne $null) {
    $rng = $targetCell.Range
} else {
    # Respaldo: si no se localizó la celda, busca en todo el documento.
    $rng = $d.Content
}

# Quita el espacio inicial antes de "companiaContabTelefono" y el espacio
# entre "companiaContabTelefono" y la llave de cierre "}".
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute(" companiaContabTelefono }", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "companiaContabTelefono}", 2)

Write-Output ("companiaContabTelefono placeholder fixed: " + $found)
